# Update diagrammes de Gantt
# Applies the changes described by the commit: adds a new "Démarrage" category
# (re-categorising the "réception Fortinet" work previously filed under
# "Documentation"), fills in the details for rows 10 & 11 of the
# "Journal de travail" sheet, adds a "-" placeholder in the "Journal de bord"
# sheet, and appends "Démarrage" to the "Paramètres" category list.

$wb = $excel.ActiveWorkbook

$wsTravail    = $wb.Worksheets.Item(1)   # "Journal de travail"
$wsBord       = $wb.Worksheets.Item(2)   # "Journal de bord"
$wsParametres = $wb.Worksheets.Item(3)   # "Paramètres"

# ---------------------------------------------------------------------
# Make sure new shared-string entries get created/appended in the same
# order as in the target workbook: "-", "Réception fortinet",
# "Réception, déballage et inventaire du matériel reçu",
# "Découverte du matériel",
# "Lecture mode d'emploi et recherche documentation supplémentaire",
# "Démarrage".
# ---------------------------------------------------------------------

# "Journal de bord" - new comment placeholder for the 3rd day
$wsBord.Range("C4").Value = "-"

# "Journal de travail" - fill in the two rows documenting the Fortinet
# delivery / unboxing work that happened on 02.02.2023
$wsTravail.Range("D10").Value = "Réception fortinet"
$wsTravail.Range("E10").Value = "Réception, déballage et inventaire du matériel reçu"

$wsTravail.Range("D11").Value = "Découverte du matériel"
$wsTravail.Range("E11").Value = "Lecture mode d'emploi et recherche documentation supplémentaire"

# New category name, first used on row 2 so the shared string is created now
$wsTravail.Range("C2").Value = "Démarrage"

# ---------------------------------------------------------------------
# Recategorize rows 2-6 from "Documentation" to the new "Démarrage" category
# ---------------------------------------------------------------------
$wsTravail.Range("C3").Value = "Démarrage"
$wsTravail.Range("C4").Value = "Démarrage"
$wsTravail.Range("C5").Value = "Démarrage"
$wsTravail.Range("C6").Value = "Démarrage"

# ---------------------------------------------------------------------
# Row 10 : 02.02.2023, 13:40 -> 14:00, Réception fortinet
# ---------------------------------------------------------------------
$wsTravail.Range("C10").Value = "Démarrage"
$wsTravail.Range("F10").Value = 0.56944444444444442
$wsTravail.Range("G10").Value = 0.58333333333333337

# ---------------------------------------------------------------------
# Row 11 : 02.02.2023, 14:00 -> 15:05, Découverte du matériel
# ---------------------------------------------------------------------
$wsTravail.Range("A11").Value = 44959
$wsTravail.Range("B11").Value = 1
$wsTravail.Range("C11").Value = "Démarrage"
$wsTravail.Range("F11").Value = 0.58333333333333337
$wsTravail.Range("G11").Value = 0.62847222222222221

# ---------------------------------------------------------------------
# "Paramètres" - append the new "Démarrage" category to the list
# (re-use the same left/center alignment as the other category cells)
# ---------------------------------------------------------------------
$wsParametres.Range("A7").Value = "Démarrage"
$wsParametres.Range("A7").HorizontalAlignment = -4131
$wsParametres.Range("A7").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Update the data-validation list on "Journal de travail"!C2:C50 so it
# references the extended Paramètres range (A2:A7) and covers the whole
# C2:C50 block (previously split as C2:C5 / C7:C50 because row 6 used to
# be excluded).
# ---------------------------------------------------------------------
$validatedRange = $wsTravail.Range("C2:C50")
$validatedRange.Validation.Delete()
$validatedRange.Validation.Add(3, 1, 1, "Paramètres!`$A`$2:`$A`$7")
$validatedRange.Validation.IgnoreBlank = $true
$validatedRange.Validation.InCellDropdown = $true
$validatedRange.Validation.ShowInput = $true
$validatedRange.Validation.ShowError = $true

# ---------------------------------------------------------------------
# Restore the selections / active cells recorded in the saved workbook
# ---------------------------------------------------------------------
$wsBord.Activate() | Out-Null
$wsBord.Range("C5").Select() | Out-Null

$wsParametres.Activate() | Out-Null
$wsParametres.Range("A8").Select() | Out-Null

$wsTravail.Activate() | Out-Null
$wsTravail.Range("E25").Select() | Out-Null
